# Update average_county_temperature (I), worst_ashp_cop (N) and best_ashp_cop (O)
# for facilities that received refreshed NOAA temperature data.
# Rows without N/O values correspond to the "not_electrifiable" option and only
# have column I updated; the paired "electrified_utilities" row (one above) also
# carries the recalculated worst/best ASHP COP values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  I = 0.4166666666666667; N = 0.988840803709428;  O = 1.017303370786517 },
    @{ Row = 3;  I = 0.4166666666666667 },
    @{ Row = 4;  I = 13.46442495126706;  N = 1.039134630545334;  O = 1.071143586266607 },
    @{ Row = 5;  I = 13.46442495126706 },
    @{ Row = 10; I = 15.36574074074072;  N = 1.046893692841948;  O = 1.079468594462633 },
    @{ Row = 11; I = 15.36574074074072 },
    @{ Row = 18; I = 14.47727272727272;  N = 1.043253580005337;  O = 1.075562343793391 },
    @{ Row = 19; I = 14.47727272727272 },
    @{ Row = 24; I = 12.41429539295394;  N = 1.034898269710531;  O = 1.066600353276151 },
    @{ Row = 25; I = 12.41429539295394 },
    @{ Row = 26; I = 16.86342592592595;  N = 1.0530876503132;    O = 1.086117960679774 },
    @{ Row = 27; I = 16.86342592592595 },
    @{ Row = 32; I = 1.791666666666668;  N = 0.9939102066179896; O = 1.022720671292561 },
    @{ Row = 33; I = 1.791666666666668 },
    @{ Row = 36; I = 15.36574074074072;  N = 1.046893692841948;  O = 1.079468594462633 },
    @{ Row = 37; I = 15.36574074074072 },
    @{ Row = 42; I = 14.47727272727272;  N = 1.043253580005337;  O = 1.075562343793391 },
    @{ Row = 43; I = 14.47727272727272 },
    @{ Row = 46; I = 17.25771604938272;  N = 1.054730517716163;  O = 1.087882159227449 },
    @{ Row = 47; I = 17.25771604938272 },
    @{ Row = 50; I = 14.47727272727272;  N = 1.043253580005337;  O = 1.075562343793391 },
    @{ Row = 51; I = 14.47727272727272 },
    @{ Row = 58; I = 14.47727272727272;  N = 1.043253580005337;  O = 1.075562343793391 },
    @{ Row = 59; I = 14.47727272727272 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("I$r").Value = $u.I
    if ($u.ContainsKey("N")) {
        $ws.Range("N$r").Value = $u.N
    }
    if ($u.ContainsKey("O")) {
        $ws.Range("O$r").Value = $u.O
    }
}
